# Regenerate the handoff report: the old source file
# "106f66bb-ac90-4da2-83c3-1c5d74c59d6f.md" was re-handed-off and is now
# tracked under a new guid, and the "240c1d6d-...md" entry (which had failed
# its handoff transform) has dropped out of the report entirely. New handoff
# packages (and timestamps) were generated for both the zh-cn and de-de
# targets of the renamed file.

$wb = $excel.ActiveWorkbook

$oldGuidFile = "106f66bb-ac90-4da2-83c3-1c5d74c59d6f.md"
$newGuidFile = "c09f9543-7afc-4c05-93c4-c9f2a3462f7c.md"

$oldHash = "618534c72b7b018ffc58208ceccaff5c5816cbc8"
$newHash = "df2a4a0026a8c5b0aa0a52ca51e11dfee4712baa"

$oldZhXlf = "106f66bb-ac90-4da2-83c3-1c5d74c59d6f.$oldHash.zh-cn.xlf"
$newZhXlf = "c09f9543-7afc-4c05-93c4-c9f2a3462f7c.$newHash.zh-cn.xlf"

$oldDeXlf = "106f66bb-ac90-4da2-83c3-1c5d74c59d6f.$oldHash.de-de.xlf"
$newDeXlf = "c09f9543-7afc-4c05-93c4-c9f2a3462f7c.$newHash.de-de.xlf"

$oldZhTime = "2016-01-19 07:09:19"
$newZhTime = "2016-01-19 07:10:07"

$oldDeTime = "2016-01-19 07:09:29"
$newDeTime = "2016-01-19 07:10:17"

$failedFile = "240c1d6d-ed74-4f0c-b90f-4eb48181bf2f.md"

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws3 = $wb.Worksheets.Item(3)

foreach ($ws in @($ws1, $ws2, $ws3)) {
    # Find the row whose "source file" column holds the failed-transform
    # entry and drop it - it is no longer part of the handoff report.
    $cell = $ws.Cells.Find($failedFile)
    if ($cell -ne $null) {
        $ws.Rows.Item($cell.Row).Delete()
    }

    # Roll the renamed source file + its regenerated handoff packages /
    # timestamps forward across every remaining cell.
    $ws.Cells.Replace($oldGuidFile, $newGuidFile)
    $ws.Cells.Replace($oldZhXlf, $newZhXlf)
    $ws.Cells.Replace($oldDeXlf, $newDeXlf)
    $ws.Cells.Replace($oldZhTime, $newZhTime)
    $ws.Cells.Replace($oldDeTime, $newDeTime)
}

# Hyperlink objects don't ride along with Replace()/row deletes, so rebuild
# each sheet's hyperlink collection from scratch against the now-correct
# grid.

$ws1.Cells.Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/98f92b1dccc0fed30b7f7661c99862cb7c13f86d/e2e/$newGuidFile", "", "", $newGuidFile)
$ws1.Hyperlinks.Add($ws1.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/98f92b1dccc0fed30b7f7661c99862cb7c13f86d/.localization-config", "", "", ".localization-config")

$ws2.Cells.Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/98f92b1dccc0fed30b7f7661c99862cb7c13f86d/e2e/$newGuidFile", "", "", $newGuidFile)
$ws2.Hyperlinks.Add($ws2.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a2544820034890bab6b7e153b5dc01cf4cf75339/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/$newZhXlf", "", "", $newZhXlf)
$ws2.Hyperlinks.Add($ws2.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/98f92b1dccc0fed30b7f7661c99862cb7c13f86d/.localization-config", "", "", ".localization-config")

$ws3.Cells.Hyperlinks.Delete()
$ws3.Hyperlinks.Add($ws3.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/98f92b1dccc0fed30b7f7661c99862cb7c13f86d/e2e/$newGuidFile", "", "", $newGuidFile)
$ws3.Hyperlinks.Add($ws3.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/77e6861c635d68370c57d4c54f24a2fc04636d4e/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/$newDeXlf", "", "", $newDeXlf)
$ws3.Hyperlinks.Add($ws3.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/98f92b1dccc0fed30b7f7661c99862cb7c13f86d/.localization-config", "", "", ".localization-config")
